$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update URL property (row 2) ---
$ws.Range("B2").Value = "http://fhir.ch/ig/ch-epr-term/ValueSet/DocumentEntry.classCode"

# --- Update Identifier property (row 3): add "use: " before OFFICIAL ---
$ws.Range("B3").Value = "id: 2.16.756.5.30.1.127.3.10.1.3 (use: OFFICIAL)"

# --- Update Version property (row 4): 1.0.0 -> 2.0.0-ballot ---
$ws.Range("B4").Value = "2.0.0-ballot"

# --- Insert a new "Contact" row (row 12), duplicating the existing Contact row (row 11) ---
$ws.Rows.Item(12).Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
